# Scheduled-runner refresh of market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# across the per-job sheets. Values below are taken straight from the
# refreshed snapshot; a handful of rows gain/lose a LeveProfit cell
# because a recipe flips between NQ-only and HQ-capable.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value2 = 74644.82000000001
$ws.Range("I70").Value2 = 2264.2
$ws.Range("J70").Value2 = 134962
$ws.Range("K70").Value2 = 6792.599999999999
$ws.Range("L70").Value2 = 404886
$ws.Range("M70").Value2 = -6522.599999999999
$ws.Range("N70").Value2 = -405426

$ws.Range("H73").Value2 = 74644.82000000001
$ws.Range("I73").Value2 = 2264.2
$ws.Range("J73").Value2 = 134962
$ws.Range("K73").Value2 = 6792.599999999999
$ws.Range("L73").Value2 = 404886
$ws.Range("M73").Value2 = -5856.599999999999
$ws.Range("N73").Value2 = -406758

$ws.Range("H76").Value2 = 7142.4287
$ws.Range("I76").Value2 = 4999
$ws.Range("J76").Value2 = 7499.6665
$ws.Range("K76").Value2 = 4999
$ws.Range("L76").Value2 = 7499.6665
$ws.Range("M76").Value2 = -4684
$ws.Range("N76").Value2 = -8129.6665

$ws.Range("H79").Value2 = 7142.4287
$ws.Range("I79").Value2 = 4999
$ws.Range("J79").Value2 = 7499.6665
$ws.Range("K79").Value2 = 4999
$ws.Range("L79").Value2 = 7499.6665
$ws.Range("M79").Value2 = -3907
$ws.Range("N79").Value2 = -9683.666499999999

$ws.Range("H93").Value2 = 29999.5
$ws.Range("J93").Value2 = 29999.5
$ws.Range("L93").Value2 = 29999.5
$ws.Range("N93").Value2 = -34991.5

$ws.Range("H106").Value2 = 15992.6
$ws.Range("I106").Value2 = 15987.667
$ws.Range("J106").Value2 = 16000
$ws.Range("K106").Value2 = 15987.667
$ws.Range("L106").Value2 = 16000
$ws.Range("M106").Value2 = -15356.667
$ws.Range("N106").Value2 = -17262

$ws.Range("H120").Value2 = 50000
$ws.Range("J120").Value2 = 50000
$ws.Range("L120").Value2 = 50000
$ws.Range("N120").Value2 = -59676

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 12138.417
$ws.Range("I32").Value2 = 9529.333000000001
$ws.Range("K32").Value2 = 9529.333000000001
$ws.Range("M32").Value2 = -9242.333000000001

$ws.Range("H97").Value2 = 1054.5
$ws.Range("J97").Value2 = 1049.5
$ws.Range("L97").Value2 = 1049.5
$ws.Range("N97").Value2 = -2041.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value2 = 3261
$ws.Range("I11").Value2 = 3133.2
$ws.Range("J11").Value2 = 3900
$ws.Range("K11").Value2 = 3133.2
$ws.Range("L11").Value2 = 3900
$ws.Range("M11").Value2 = -2993.2
$ws.Range("N11").Value2 = -4180

$ws.Range("H94").Value2 = 2159.1
$ws.Range("I94").Value2 = 2065.6667
$ws.Range("K94").Value2 = 2065.6667
$ws.Range("M94").Value2 = -1614.6667

$ws.Range("H105").Value2 = 1522.9445
$ws.Range("I105").Value2 = 1150.875
$ws.Range("J105").Value2 = 4499.5
$ws.Range("K105").Value2 = 1150.875
$ws.Range("L105").Value2 = 4499.5
$ws.Range("M105").Value2 = 596.125
$ws.Range("N105").Value2 = -7993.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 7995.75
$ws.Range("I4").Value2 = 0
$ws.Range("K4").Value2 = 0
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value2 = 3259.2778
$ws.Range("I31").Value2 = 1741
$ws.Range("K31").Value2 = 1741
$ws.Range("M31").Value2 = -1446

$ws.Range("H34").Value2 = 3259.2778
$ws.Range("I34").Value2 = 1741
$ws.Range("K34").Value2 = 1741
$ws.Range("M34").Value2 = -1539

$ws.Range("H43").Value2 = 16942.625
$ws.Range("J43").Value2 = 17684.428
$ws.Range("L43").Value2 = 17684.428
$ws.Range("N43").Value2 = -18052.428

$ws.Range("H95").Value2 = 28032.111
$ws.Range("J95").Value2 = 28032.111
$ws.Range("L95").Value2 = 28032.111
$ws.Range("N95").Value2 = -33524.111

$ws.Range("H101").Value2 = 16942.625
$ws.Range("J101").Value2 = 17684.428
$ws.Range("L101").Value2 = 17684.428
$ws.Range("N101").Value2 = -24174.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 1250614.1
$ws.Range("I4").Value2 = 2125466.5
$ws.Range("K4").Value2 = 6376399.5
$ws.Range("M4").Value2 = -6376287.5

$ws.Range("H62").Value2 = 4099.875
$ws.Range("I62").Value2 = 8500
$ws.Range("J62").Value2 = 3471.2856
$ws.Range("K62").Value2 = 25500
$ws.Range("L62").Value2 = 10413.8568
$ws.Range("M62").Value2 = -24814
$ws.Range("N62").Value2 = -11785.8568

$ws.Range("H65").Value2 = 4099.875
$ws.Range("I65").Value2 = 8500
$ws.Range("J65").Value2 = 3471.2856
$ws.Range("K65").Value2 = 76500
$ws.Range("L65").Value2 = 31241.5704
$ws.Range("M65").Value2 = -73068
$ws.Range("N65").Value2 = -38105.5704

$ws.Range("H96").Value2 = 0
$ws.Range("J96").Value2 = 0
$ws.Range("L96").Value2 = 0
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 4600200
$ws.Range("I11").Value2 = 6666999.5
$ws.Range("J11").Value2 = 3714428.5
$ws.Range("K11").Value2 = 6666999.5
$ws.Range("L11").Value2 = 3714428.5
$ws.Range("M11").Value2 = -6666860.5
$ws.Range("N11").Value2 = -3714706.5

$ws.Range("H97").Value2 = 566.5
$ws.Range("I97").Value2 = 620.8
$ws.Range("K97").Value2 = 620.8
$ws.Range("M97").Value2 = -124.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value2 = 19000
$ws.Range("J3").Value2 = 19000
$ws.Range("L3").Value2 = 19000
$ws.Range("N3").Value2 = -19224

$ws.Range("H15").Value2 = 19000
$ws.Range("J15").Value2 = 19000
$ws.Range("L15").Value2 = 19000
$ws.Range("N15").Value2 = -19340

$ws.Range("H20").Value2 = 29420
$ws.Range("I20").Value2 = 29420
$ws.Range("K20").Value2 = 29420
$ws.Range("M20").Value2 = -29194

$ws.Range("H21").Value2 = 12400
$ws.Range("J21").Value2 = 12400
$ws.Range("L21").Value2 = 12400
$ws.Range("N21").Value2 = -12748

$ws.Range("H68").Value2 = 3067
$ws.Range("J68").Value2 = 3500.5
$ws.Range("L68").Value2 = 3500.5
$ws.Range("N68").Value2 = -4998.5

$ws.Range("H71").Value2 = 3067
$ws.Range("J71").Value2 = 3500.5
$ws.Range("L71").Value2 = 17502.5
$ws.Range("N71").Value2 = -24990.5

$ws.Range("H93").Value2 = 1426.25
$ws.Range("I93").Value2 = 1185.2
$ws.Range("K93").Value2 = 1185.2
$ws.Range("M93").Value2 = 62.79999999999995

$ws.Range("H122").Value2 = 6454.6665
$ws.Range("I122").Value2 = 7045.6
$ws.Range("K122").Value2 = 21136.8
$ws.Range("M122").Value2 = -18686.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value2 = 15999
$ws.Range("I58").Value2 = 15999
$ws.Range("K58").Value2 = 15999
$ws.Range("M58").Value2 = -15691
